$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value2 = '246.18'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value2 = '0.94%'
# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value2 = '29.52'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value2 = '-1.35%'
# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value2 = '5.160'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value2 = '0.01%'
# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value2 = '0.05797'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value2 = '2.15%'
# Row 6
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value2 = '1.71%'
# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value2 = '3.206'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value2 = '5.61%'
# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value2 = '0.8498'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value2 = '0.08%'
# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value2 = '0.8638'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value2 = '-0.30%'
# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value2 = '0.1386'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value2 = '2.80%'
# Row 11
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value2 = '2.59%'
# Row 12
$ws.Range('B12').Value2 = 'BitrueCoin'
$ws.Range('C12').Value2 = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value2 = '0.03200'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value2 = '10.59%'
# Row 13
$ws.Range('B13').Value2 = 'BitMartToken'
$ws.Range('C13').Value2 = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value2 = '0.09380'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value2 = '-0.01%'
# Row 14
$ws.Range('B14').Value2 = 'BitForexToken'
$ws.Range('C14').Value2 = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value2 = '0.001541'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value2 = '0.90%'
# Row 15
$ws.Range('B15').Value2 = 'One'
$ws.Range('C15').Value2 = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value2 = '0.0006026'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value2 = '0.30%'
# Row 16
$ws.Range('B16').Value2 = 'TigerCash'
$ws.Range('C16').Value2 = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value2 = '0.006121'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value2 = '0.12%'
# Row 17
$ws.Range('B17').Value2 = 'LEO'
$ws.Range('C17').Value2 = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value2 = '3.488'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value2 = '-0.64%'
# Row 18
$ws.Range('B18').Value2 = 'BTSEToken'
$ws.Range('C18').Value2 = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value2 = '2.221'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value2 = '-1.07%'
# Row 19
$ws.Range('B19').Value2 = 'BitpandaEcosystemToken'
$ws.Range('C19').Value2 = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value2 = '0.3197'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value2 = '1.56%'
# Row 20
$ws.Range('B20').Value2 = 'LiechtensteinCryptoassetsExchange'
$ws.Range('C20').Value2 = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value2 = '0.03366'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value2 = '0.44%'
# Row 21
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value2 = '-1.60%'
# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value2 = '3.489'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value2 = '-3.83%'
# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value2 = '0.04150'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value2 = '-0.65%'
# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value2 = '0.001228'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value2 = '1.43%'
# Row 26
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value2 = '-6.73%'
# Row 28
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value2 = '3.85%'
# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value2 = '0.03749'
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value2 = '-1.15%'
# Row 41
$ws.Range('B41').Value2 = 'KickToken'
$ws.Range('C41').Value2 = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value2 = '0.005786'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value2 = '-0.38%'
# Row 42
$ws.Range('B42').Value2 = 'BKEXToken'
$ws.Range('C42').Value2 = 'https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value2 = '0.1071'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value2 = '1.36%'
# Row 43
$ws.Range('B43').Value2 = 'CEJI'
$ws.Range('C43').Value2 = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value2 = '0.002199'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value2 = '-3.94%'
# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value2 = '0.009189'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value2 = '-1.02%'
# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value2 = '0.00005299'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value2 = '3.87%'
# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value2 = '0.00000000750'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value2 = '-0.03%'
# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value2 = '0.05797'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value2 = '-35.57%'
# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value2 = '0.002175'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value2 = '-21.36%'
# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value2 = '0.00002099'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value2 = '-0.03%'
# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value2 = '0.0001999'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value2 = '-0.03%'
